$d = $word.ActiveDocument

# Make sure revisions are tracked so the replacements are recorded as
# w:ins / w:del pairs (the document already has trackRevisions set, but
# be explicit).
$d.TrackRevisions = $true

# Word keeps a single hidden "_GoBack" bookmark that marks the location
# of the most recent edit. Remove the one that currently lives inside
# the "goal of the Agile Link SDK" paragraph; it will be recreated below
# around the final edit we make, exactly as Word would do as edits are
# applied in document order.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Replace every whole-word occurrence of "Agile Link" (with the space)
# with "AMAP" as a tracked change.
$d.Content.Find.Execute("Agile Link", $true, $false, $false, $false, $false, $true, 1, $false, "AMAP", 2)

# Replace the one-word spelling "AgileLink" (no space) with "AMAP" too.
$d.Content.Find.Execute("AgileLink", $true, $false, $false, $false, $false, $true, 1, $false, "AMAP", 2)

# Re-create the "_GoBack" bookmark around the last edit that was made
# (the AgileLink -> AMAP replacement), matching Word's behaviour of
# tracking the most recent editing location.
$found = $d.Content.Find.Execute("AgileLink", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
